$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ALC.Range("H2").Value = 1600.2333
$ws_ALC.Range("I2").Value = 363
$ws_ALC.Range("J2").Value = 1847.68
$ws_ALC.Range("K2").Value = 363
$ws_ALC.Range("L2").Value = 1847.68
$ws_ALC.Range("M2").Value = -250
$ws_ALC.Range("N2").Value = -2073.68

$ws_ALC.Range("H11").Value = 10290.533
$ws_ALC.Range("I11").Value = 10290.533
$ws_ALC.Range("K11").Value = 10290.533
$ws_ALC.Range("M11").Value = -10150.533

$ws_ALC.Range("H38").Value = 2434.5
$ws_ALC.Range("J38").Value = 4187.25
$ws_ALC.Range("L38").Value = 12561.75
$ws_ALC.Range("N38").Value = -13305.75

$ws_ALC.Range("H40").Value = 2507.8333
$ws_ALC.Range("I40").Value = 1800.5
$ws_ALC.Range("J40").Value = 2861.5
$ws_ALC.Range("K40").Value = 1800.5
$ws_ALC.Range("L40").Value = 2861.5
$ws_ALC.Range("M40").Value = -1625.5
$ws_ALC.Range("N40").Value = -3211.5

$ws_ALC.Range("H53").Value = 45455020
$ws_ALC.Range("I53").Value = 451.54544
$ws_ALC.Range("J53").Value = 90909590
$ws_ALC.Range("K53").Value = 451.54544
$ws_ALC.Range("L53").Value = 90909590
$ws_ALC.Range("M53").Value = 185.45456
$ws_ALC.Range("N53").Value = -90910864

$ws_ALC.Range("H58").Value = 811.75

$ws_ALC.Range("H64").Value = 7699.5713
$ws_ALC.Range("I64").Value = 7013
$ws_ALC.Range("K64").Value = 7013
$ws_ALC.Range("M64").Value = -6765

$ws_ALC.Range("H67").Value = 7699.5713
$ws_ALC.Range("I67").Value = 7013
$ws_ALC.Range("K67").Value = 7013
$ws_ALC.Range("M67").Value = -6155

$ws_ALC.Range("H132").Value = 1993.9215
$ws_ALC.Range("I132").Value = 1886.8043
$ws_ALC.Range("K132").Value = 5660.4129
$ws_ALC.Range("M132").Value = -3130.4129

$ws_ALC.Range("H138").Value = 5027.3906
$ws_ALC.Range("I138").Value = 7921.227
$ws_ALC.Range("J138").Value = 3511.5715
$ws_ALC.Range("K138").Value = 23763.681
$ws_ALC.Range("L138").Value = 10534.7145
$ws_ALC.Range("M138").Value = -18623.681
$ws_ALC.Range("N138").Value = -20814.7145

$ws_ALC.Range("H141").Value = 2499.75
$ws_ALC.Range("I141").Value = 2499.75
$ws_ALC.Range("K141").Value = 7499.25
$ws_ALC.Range("M141").Value = -2319.25

$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_ARM.Range("H61").Value = 1595447.2
$ws_ARM.Range("J61").Value = 3100775.8
$ws_ARM.Range("L61").Value = 3100775.8
$ws_ARM.Range("N61").Value = -3101199.8

$ws_ARM.Range("H102").Value = 0
$ws_ARM.Range("I102").Value = 0
$ws_ARM.Range("K102").Value = 0
$ws_ARM.Range("M102").ClearContents()

$ws_ARM.Range("H105").Value = 124999.5
$ws_ARM.Range("J105").Value = 124999.5
$ws_ARM.Range("L105").Value = 124999.5
$ws_ARM.Range("N105").Value = -131987.5

$ws_ARM.Range("H110").Value = 704.7273
$ws_ARM.Range("I110").Value = 704.7273
$ws_ARM.Range("K110").Value = 704.7273
$ws_ARM.Range("M110").Value = 1340.2727

$ws_ARM.Range("H136").Value = 1595447.2
$ws_ARM.Range("J136").Value = 3100775.8
$ws_ARM.Range("L136").Value = 9302327.399999999
$ws_ARM.Range("N136").Value = -9307427.399999999

$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_BSM.Range("H20").Value = 873.23334
$ws_BSM.Range("I20").Value = 743
$ws_BSM.Range("J20").Value = 1231.375
$ws_BSM.Range("K20").Value = 743
$ws_BSM.Range("L20").Value = 1231.375
$ws_BSM.Range("M20").Value = -496
$ws_BSM.Range("N20").Value = -1725.375

$ws_BSM.Range("H99").Value = 5673.9585
$ws_BSM.Range("I99").Value = 10831.777
$ws_BSM.Range("J99").Value = 2579.2666
$ws_BSM.Range("K99").Value = 10831.777
$ws_BSM.Range("L99").Value = 2579.2666
$ws_BSM.Range("M99").Value = -9333.777
$ws_BSM.Range("N99").Value = -5575.2666

$ws_BSM.Range("H105").Value = 8466.210999999999
$ws_BSM.Range("J105").Value = 4245.6665
$ws_BSM.Range("L105").Value = 4245.6665
$ws_BSM.Range("N105").Value = -7739.6665

$ws_BSM.Range("H107").Value = 19383.941
$ws_BSM.Range("I107").Value = 21823.857
$ws_BSM.Range("K107").Value = 21823.857
$ws_BSM.Range("M107").Value = -19903.857

$ws_BSM.Range("H134").Value = 17649234
$ws_BSM.Range("I134").Value = 1631.2439
$ws_BSM.Range("K134").Value = 4893.7317
$ws_BSM.Range("M134").Value = -2358.7317

$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CRP.Range("H31").Value = 3003.554
$ws_CRP.Range("I31").Value = 1814.5454
$ws_CRP.Range("J31").Value = 3245.7593
$ws_CRP.Range("K31").Value = 1814.5454
$ws_CRP.Range("L31").Value = 3245.7593
$ws_CRP.Range("M31").Value = -1519.5454
$ws_CRP.Range("N31").Value = -3835.7593

$ws_CRP.Range("H34").Value = 3003.554
$ws_CRP.Range("I34").Value = 1814.5454
$ws_CRP.Range("J34").Value = 3245.7593
$ws_CRP.Range("K34").Value = 1814.5454
$ws_CRP.Range("L34").Value = 3245.7593
$ws_CRP.Range("M34").Value = -1612.5454
$ws_CRP.Range("N34").Value = -3649.7593

$ws_CRP.Range("H132").Value = 15154711
$ws_CRP.Range("J132").Value = 23811944
$ws_CRP.Range("L132").Value = 71435832
$ws_CRP.Range("N132").Value = -71440892

$ws_CRP.Range("H134").Value = 2081.56
$ws_CRP.Range("I134").Value = 1813
$ws_CRP.Range("K134").Value = 5439
$ws_CRP.Range("M134").Value = -2904

$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_CUL.Range("H92").Value = 594
$ws_CUL.Range("I92").Value = 538.25
$ws_CUL.Range("J92").Value = 649.75
$ws_CUL.Range("K92").Value = 1614.75
$ws_CUL.Range("L92").Value = 1949.25
$ws_CUL.Range("M92").Value = -366.75
$ws_CUL.Range("N92").Value = -4445.25

$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_GSM.Range("H80").Value = 323192.6
$ws_GSM.Range("I80").Value = 433665.66
$ws_GSM.Range("K80").Value = 433665.66
$ws_GSM.Range("M80").Value = -432667.66

$ws_GSM.Range("H83").Value = 323192.6
$ws_GSM.Range("I83").Value = 433665.66
$ws_GSM.Range("K83").Value = 2168328.3
$ws_GSM.Range("M83").Value = -2163336.3

$ws_GSM.Range("H137").Value = 179990
$ws_GSM.Range("J137").Value = 179990
$ws_GSM.Range("L137").Value = 179990
$ws_GSM.Range("N137").Value = -190190

$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_LTW.Range("H122").Value = 5431.3438
$ws_LTW.Range("I122").Value = 3990.8333
$ws_LTW.Range("J122").Value = 6295.65
$ws_LTW.Range("K122").Value = 11972.4999
$ws_LTW.Range("L122").Value = 18886.95
$ws_LTW.Range("M122").Value = -9522.499899999999
$ws_LTW.Range("N122").Value = -23786.95

$ws_LTW.Range("H132").Value = 35333.332
$ws_LTW.Range("I132").Value = 3000
$ws_LTW.Range("K132").Value = 9000
$ws_LTW.Range("M132").Value = -6470

$ws_LTW.Range("H136").Value = 5251.706
$ws_LTW.Range("J136").Value = 6252.231
$ws_LTW.Range("L136").Value = 18756.693
$ws_LTW.Range("N136").Value = -23856.693

$ws_WVR = $wb.Worksheets.Item("WVR")
$ws_WVR.Range("H3").Value = 252775
$ws_WVR.Range("I3").Value = 0
$ws_WVR.Range("J3").Value = 252775
$ws_WVR.Range("K3").Value = 0
$ws_WVR.Range("L3").Value = 252775
$ws_WVR.Range("N3").Value = -253003
$ws_WVR.Range("M3").ClearContents()

$ws_WVR.Range("H49").Value = 173352
$ws_WVR.Range("I49").Value = 173352
$ws_WVR.Range("K49").Value = 173352
$ws_WVR.Range("M49").Value = -173122

$ws_WVR.Range("H81").Value = 3240968.2
$ws_WVR.Range("I81").Value = 4536679
$ws_WVR.Range("J81").Value = 1692.25
$ws_WVR.Range("K81").Value = 9073358
$ws_WVR.Range("L81").Value = 3384.5
$ws_WVR.Range("M81").Value = -9072297
$ws_WVR.Range("N81").Value = -5506.5

$ws_WVR.Range("H84").Value = 3240968.2
$ws_WVR.Range("I84").Value = 4536679
$ws_WVR.Range("J84").Value = 1692.25
$ws_WVR.Range("K84").Value = 45366790
$ws_WVR.Range("L84").Value = 16922.5
$ws_WVR.Range("M84").Value = -45361486
$ws_WVR.Range("N84").Value = -27530.5

$ws_WVR.Range("H126").Value = 2432.65
$ws_WVR.Range("I126").Value = 2229.5881
$ws_WVR.Range("K126").Value = 6688.7643
$ws_WVR.Range("M126").Value = -4218.7643

$ws_WVR.Range("H132").Value = 40461
$ws_WVR.Range("I132").Value = 54505.316
$ws_WVR.Range("K132").Value = 163515.948
$ws_WVR.Range("M132").Value = -160985.948

$ws_WVR.Range("H136").Value = 48060.684
$ws_WVR.Range("I136").Value = 112560.664
$ws_WVR.Range("J136").Value = 3406.8462
$ws_WVR.Range("K136").Value = 337681.992
$ws_WVR.Range("L136").Value = 10220.5386
$ws_WVR.Range("M136").Value = -335131.992
$ws_WVR.Range("N136").Value = -15320.5386
